$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.175720741260731
$ws.Range("C2").Value = 1.175720741260731
$ws.Range("D2").Value = 2.00046158807684
$ws.Range("E2").Value = 1.414376748987638
$ws.Range("F2").Value = 0.8158993420687196

$ws.Range("B3").Value = 1.179154357706104
$ws.Range("C3").Value = 1.179154357706104
$ws.Range("D3").Value = 1.839172490612837
$ws.Range("E3").Value = 1.356160938315522
$ws.Range("F3").Value = 0.7061377669450614

$ws.Range("B4").Value = 0.9948462645758421
$ws.Range("C4").Value = 0.9948462645758421
$ws.Range("D4").Value = 1.24218013291981
$ws.Range("E4").Value = 1.114531351250296
$ws.Range("F4").Value = 0.5504118924361679

$ws.Range("B5").Value = 1.219879926150666
$ws.Range("C5").Value = 1.219879926150666
$ws.Range("D5").Value = 1.817306308435718
$ws.Range("E5").Value = 1.348075038132417
$ws.Range("F5").Value = 0.8114176165333894
